$d = $word.ActiveDocument

$b64 = "PD94bWwgdmVyc2lvbj0iMS4wIiBlbmNvZGluZz0iVVRGLTgiIHN0YW5kYWxvbmU9InllcyI/Pjxwa2c6cGFja2FnZSB4bWxuczpwa2c9Imh0dHA6Ly9zY2hlbWFzLm1pY3Jvc29mdC5jb20vb2ZmaWNlLzIwMDYveG1sUGFja2FnZSI+PHBrZzpwYXJ0IHBrZzpuYW1lPSIvd29yZC9kb2N1bWVudC54bWwiIHBrZzpjb250ZW50VHlwZT0iYXBwbGljYXRpb24vdm5kLm9wZW54bWxmb3JtYXRzLW9mZmljZWRvY3VtZW50LndvcmRwcm9jZXNzaW5nbWwuZG9jdW1lbnQubWFpbit4bWwiIHBrZzpwYWRkaW5nPSI1MTIiPjxwa2c6eG1sRGF0YT48dzpkb2N1bWVudCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OmJvZHk+PHc6cC8+Cjx3OnAvPgo8dzpwLz4KPHc6cC8+Cjx3OnAvPgo8dzpwLz4KPHc6cD4KICA8dzpwUHI+CiAgICA8dzpwU3R5bGUgdzp2YWw9IkhlYWRpbmcxIi8+CiAgPC93OnBQcj4KICA8dzpyPgogICAgPHc6bGFzdFJlbmRlcmVkUGFnZUJyZWFrLz4KICAgIDx3OnQ+VG9waWMgMjogRXhwYW5kZWQgVXNlIENhc2U8L3c6dD4KICA8L3c6cj4KPC93OnA+Cjx3OnAvPgo8dzpwPgogIDx3OnBQcj4KICAgIDx3OnBTdHlsZSB3OnZhbD0iSGVhZGluZzIiLz4KICA8L3c6cFByPgogIDx3OnI+CiAgICA8dzp0PlNsaWRlIDE6PC93OnQ+CiAgPC93OnI+CjwvdzpwPgo8dzpwPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPkhlbGxvLCBteSBuYW1lIGlzIDwvdzp0PgogIDwvdzpyPgogIDx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxTdGFydCIvPgogIDx3OnI+CiAgICA8dzp0PkthbmFnYTwvdzp0PgogIDwvdzpyPgogIDx3OnByb29mRXJyIHc6dHlwZT0ic3BlbGxFbmQiLz4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gTWFuaWthbmRhbi4gVG9kYXksIEkgYW0gZ29pbmcgdG8gdGFsayBhYm91dCBleHBhbmRlZCB1c2UgY2FzZXMgYW5kIGhvdyB0aGV5IHdlcmUgdXNlZCBpbiB0aGUgc29mdHdhcmUgZGV2ZWxvcG1lbnQgZXhlcmNpc2UuIDwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzpyUHI+CiAgICAgIDx3OmNvbG9yIHc6dmFsPSJGRjAwMDAiLz4KICAgIDwvdzpyUHI+CiAgICA8dzp0PkNsaWNrPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQ+Ljwvdzp0PgogIDwvdzpyPgo8L3c6cD4KPHc6cC8+Cjx3OnA+CiAgPHc6cFByPgogICAgPHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMiIvPgogIDwvdzpwUHI+CiAgPHc6cj4KICAgIDx3OnQ+U2xpZGUgMjo8L3c6dD4KICA8L3c6cj4KPC93OnA+Cjx3OnA+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dD4uIFNvLCB3aGF0IGFyZSBleHBhbmRlZCB1c2UgY2FzZXM8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4/IEV4cGFuZGVkIFVzZSBDYXNlIGlzIGEgZGV0YWlsZWQgZGVzY3JpcHRpb24gb2YgdGhlIHByb2Nlc3NlcyB1c2VkIHRvIGNvbXBsZXRlIHZhcmlvdXMgc3lzdGVtIGZ1bmN0aW9ucy4gPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4uIFdl4oCZbGwgYmUgY292ZXJpbmcgb24gaG93IGV4cGFuZGVkIHVzZSBjYXNlcyB3ZXJlIHVzZWQgdG8gaWRlbnRpZnkgcmVxdWlyZW1lbnRzIGFuZCBvcmdhbmlzZSB0aGUgbWFuYWdlbWVudCBvZiB0aGUgc29mdHdhcmUgZGV2ZWxvcG1lbnQgZXhlcmNpc2UuIDwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzpyUHI+CiAgICAgIDx3OmNvbG9yIHc6dmFsPSJGRjAwMDAiLz4KICAgIDwvdzpyUHI+CiAgICA8dzp0PkNsaWNrPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+LiBXZeKAmWxsIGFsc28gYmUgY292ZXJpbmcgb24gaG93IGl0IHByb3ZpZGVkIG1vcmUga25vd2xlZGdlIHRoYW4gdGhhdCBvZiBhbiBleGlzdGluZyB1c2UgY2FzZSBkaWFncmFtLiA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6clByPgogICAgICA8dzpjb2xvciB3OnZhbD0iRkYwMDAwIi8+CiAgICA8L3c6clByPgogICAgPHc6dD5DbGljazwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0Pi48L3c6dD4KICA8L3c6cj4KPC93OnA+Cjx3OnAvPgo8dzpwPgogIDx3OnBQcj4KICAgIDx3OnBTdHlsZSB3OnZhbD0iSGVhZGluZzIiLz4KICA8L3c6cFByPgogIDx3OnI+CiAgICA8dzp0PlNsaWRlIDM6PC93OnQ+CiAgPC93OnI+CjwvdzpwPgo8dzpwPgogIDx3OnI+CiAgICA8dzpyUHI+CiAgICAgIDx3OmNvbG9yIHc6dmFsPSJGRjAwMDAiLz4KICAgIDwvdzpyUHI+CiAgICA8dzp0PkNsaWNrPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+LiBOb3csIGxldOKAmXMgdGFsayBhYm91dCBob3cgZXhwYW5kZWQgdXNlIGNhc2UgaGVscGVkIHRoZSB0ZWFtIGlkZW50aWZ5IHRoZSByZXF1aXJlbWVudHMuIFRoZSB0ZWFtIHVzZWQgdGhlIGZlYXR1cmVzIGZyb20gdGhlIDwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0PnVzZXIgc3Rvcmllczwvdzp0PgogIDwvdzpyPgogIDx3OmJvb2ttYXJrU3RhcnQgdzppZD0iMSIgdzpuYW1lPSJfR29CYWNrIi8+CiAgPHc6Ym9va21hcmtFbmQgdzppZD0iMSIvPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiB0byBkZXJpdmUgdGhlIHNvZnR3YXJl4oCZcyBmdW5jdGlvbnMgYW5kIHByb2Nlc3Nlcy4gPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4uIEZyb20gdGhlc2UgZnVuY3Rpb25zLCB0aGUgdGVhbSB3YXMgYWJsZSB0byA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dD5jYXRlZ29yaXNlPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IGVhY2ggZnVuY3Rpb24gaW50byBwcmltYXJ5IGFuZCBzZWNvbmRhcnkgZnVuY3Rpb25zLiA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6clByPgogICAgICA8dzpjb2xvciB3OnZhbD0iRkYwMDAwIi8+CiAgICA8L3c6clByPgogICAgPHc6dD5DbGljazwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPi4gUHJpbWFyeSBmdW5jdGlvbnMgYXJlIHRob3NlIHdoaWNoIGFyZSB2ZXJ5IGVzc2VudGlhbCB0byB0aGUgc29mdHdhcmUgcHJvZHVjdC4gU2Vjb25kYXJ5IGZlYXR1cmVzIGFyZSB0aG9zZSB3aGljaCB3aWxsIGhlbHAgdGhlIHVzZXIgdG8gcGVyZm9ybSBleHRyYSBmZWF0dXJlcy4gPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+VGhlIHNvZnR3YXJlIHdpbGwgYmUgYWJsZSB0byB3b3JrIGV2ZW4gd2l0aG91dCB0aGVzZSBmdW5jdGlvbnMgYnV0IGl04oCZcyBiZXR0ZXIgdG8gaGF2ZSB0aGVtIGluY2x1ZGVkLiA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5Ob3csIHRoYXQgdGhlIHRlYW0gaGFkIGRlcml2ZWQgdGhlIHR3byB0eXBlcyBvZiBmdW5jdGlvbnMsIGl0IHdhcyB0aW1lIHRvIGNvbWUgdXAgd2l0aCBhIHBsYW4gb24gaG93IHRoZXNlIGZ1bmN0aW9ucyBhcmUgZ29pbmcgdG8gYmUgY2FycmllZCBvdXQgYW5kIHdoYXQgbmVjZXNzYXJ5IHN0ZXBzIHdpbGwgaGF2ZSB0byBiZSB0YWtlbiA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dD5ieTwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiB0aGUgdXNlciB0byBjb21wbGV0ZSBhIHNwZWNpZmljIGZ1bmN0aW9uLjwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6clByPgogICAgICA8dzpjb2xvciB3OnZhbD0iRkYwMDAwIi8+CiAgICA8L3c6clByPgogICAgPHc6dD5DbGljazwvdzp0PgogIDwvdzpyPgo8L3c6cD4KPHc6cC8+Cjx3OnA+CiAgPHc6cFByPgogICAgPHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMiIvPgogIDwvdzpwUHI+CiAgPHc6cj4KICAgIDx3OnQ+U2xpZGUgNDo8L3c6dD4KICA8L3c6cj4KPC93OnA+Cjx3OnA+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4uIFdlIHdpbGwgbm93IHNlZSBob3cgdGhlIHVzZSBjYXNlIGRpYWdyYW0gYmVjYW1lIGluZmVyaW9yIHRvIHRoZSBleHBhbmRlZCB1c2UgY2FzZS4gPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dD4uIFRoaXMgdXNlIGNhc2UgZGlhZ3JhbSBwcm92aWRlZCBhIGNsZWFyIHZpc3VhbCBmbG93IG9mIHRoZSBzdGVwcyB0YWtlbiBieSB0aGUgdXNlciB0byBjb21wbGV0ZSB0aGUgcmVnaXN0cmF0aW9uIG9mIGEgbmV3IG1lbWJlciBidXQgYXMgeW91IGNhbiBzZWUsIGl0IHdhcyB1bmFibGUgdG8gZGVzY3JpYmUgc29tZSBzdGVwcyBpbiBtb3JlPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IGRldGFpbDwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0Pi48L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gVGhpcyB3aWxsIDwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPm9ubHkgYWxsb3cgdGhlIHRlYW0gdG8gdW5kZXJzdGFuZCB0aGUgZmxvdyBidXQgbm90IGNhcHR1cmUgdGhlIGNvbXBsZXRlIHBpY3R1cmUuIDwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzpyUHI+CiAgICAgIDx3OmNvbG9yIHc6dmFsPSJGRjAwMDAiLz4KICAgIDwvdzpyUHI+CiAgICA8dzp0PkNsaWNrPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+LiBJZiBtb3JlIGZlYXR1cmVzIHdlcmUgdG8gYmUgYWRkZWQsIHRoZSBkaWFncmFtIGJlY29tZXMgbW9yZSBjbHV0dGVyZWQgYW5kIHdpbGwgb25seSBnZXQgYmlnZ2VyLiA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6clByPgogICAgICA8dzpjb2xvciB3OnZhbD0iRkYwMDAwIi8+CiAgICA8L3c6clByPgogICAgPHc6dD5DbGljazwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPi4gQXMgeW91IGNhbiBzZWUsIHRoZSBleHBhbmRlZCB1c2UgY2FzZSBkaXNwbGF5cyBhIGxvdCBtb3JlIGluZm9ybWF0aW9uIG9uIGEgc2luZ2xlIGZ1bmN0aW9uIGFuZCB0aGUgdHlwaWNhbCBjb3Vyc2Ugb2YgZXZlbnRzIGRpc3BsYXlzIGJvdGggd2hhdCB0aGUgdXNlciBkb2VzIGFuZCB3aGF0IHRoZSBzeXN0ZW0gZG9lcyBpbiByZXNwb25zZS4gQWZ0ZXIgdGhpcywgdGhlIHRlYW0gd2FzIGFibGUgdG8gZ2FpbiBtb3JlIGtub3dsZWRnZSBvbiB0aGUgc2FtZSBmdW5jdGlvbiBhbmQgdW5kZXJzdGFuZCB0aGUgY29uY2VwdCBpbiBtb3JlIGRlcHRoLiA8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6clByPgogICAgICA8dzpjb2xvciB3OnZhbD0iRkYwMDAwIi8+CiAgICA8L3c6clByPgogICAgPHc6dD5DbGljazwvdzp0PgogIDwvdzpyPgo8L3c6cD4KPHc6cC8+Cjx3OnA+CiAgPHc6cFByPgogICAgPHc6cFN0eWxlIHc6dmFsPSJIZWFkaW5nMiIvPgogIDwvdzpwUHI+CiAgPHc6cj4KICAgIDx3OnQ+U2xpZGUgNTo8L3c6dD4KICA8L3c6cj4KPC93OnA+Cjx3OnA+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dD4uIE5vdywgbGV04oCZcyBzZWUgaG93IHRoaXMgd2FzIHVzZWQgdG8gb3JnYW5pc2UgdGhlIG1hbmFnZW1lbnQgb2YgdGhlIHNvZnR3YXJlIGRldmVsb3BtZW50IGV4ZXJjaXNlLiBUaGUgZGV0YWlsZWQgZGVzY3JpcHRpb24gZnJvbSB0aGUgZXhwYW5kZWQgdXNlIGNhc2UgaGVscGVkIHRoZSB0ZWFtIHRvIHNvcnQgb3V0IHRoZTwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiB3ZWlnaHQgb2YgdGhlPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IHRhc2tzIGFuZCBhc3NpZ24gdGhvc2UgdG8gZWFjaCBtZW1iZXIgZXF1YWxseS48L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dD4uIFRoZSB0eXBpY2FsIGNvdXJzZSBvZiBldmVudHM8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gZnJvbSB0aGUgZXhwYW5kZWQgdXNlIGNhc2UgPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+ZW5hYmxlZCB0aGUgdGVhbSB0byB1bmRlcnN0YW5kIHRoZSBmbG93IG9mIHRoZSBwcm9jZXNzZXMgaW4gbW9yZSBkZXRhaWwgYW5kIHN0cnVjdHVyZSB0aGUgcXVlcnkgYXBwcm9wcmlhdGVseS4gPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnJQcj4KICAgICAgPHc6Y29sb3Igdzp2YWw9IkZGMDAwMCIvPgogICAgPC93OnJQcj4KICAgIDx3OnQ+Q2xpY2s8L3c6dD4KICA8L3c6cj4KICA8dzpyPgogICAgPHc6dD4uPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IDwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzpyUHI+CiAgICAgIDx3OmNvbG9yIHc6dmFsPSJGRjAwMDAiLz4KICAgIDwvdzpyUHI+CiAgICA8dzp0PkNsaWNrPC93OnQ+CiAgPC93OnI+CiAgPHc6cj4KICAgIDx3OnQ+Ljwvdzp0PgogIDwvdzpyPgogIDx3OnI+CiAgICA8dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiA8L3c6dD4KICA8L3c6cj4KPC93OnA+Cjwvdzpib2R5Pjwvdzpkb2N1bWVudD48L3BrZzp4bWxEYXRhPjwvcGtnOnBhcnQ+PC9wa2c6cGFja2FnZT4="
$bytes = [Convert]::FromBase64String($b64)
$xml = [System.Text.Encoding]::UTF8.GetString($bytes)

$p = $d.Paragraphs(28)
$r = $d.Range($p.Range.Start, $d.Paragraphs(31).Range.End)
$r.InsertXML($xml)
